# Jos Buttler sheet: shift the per-match stat rows down by one
# (row2's old stats move to row3, row3's to row4, ... row5's to row6)
# and give row2 a new set of match stats.
#
# The sheet stores runs/balls/fours/sixes as text-typed numbers
# (numberStoredAsText), so we must write these as text rather than
# letting Excel auto-convert the numeric-looking strings to real
# numbers. We stage each value as a `="n"` formula (guaranteeing a
# text result) and then flatten the whole block to values in one
# Copy/PasteSpecial so the final cells hold plain text values with no
# residual formulas and no stray number-format/style changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Formula = '="24"'
$ws.Range("D2").Formula = '="25"'
$ws.Range("E2").Formula = '="1"'
$ws.Range("F2").Formula = '="1"'

$ws.Range("C3").Formula = '="22"'
$ws.Range("D3").Formula = '="11"'
$ws.Range("E3").Formula = '="1"'
$ws.Range("F3").Formula = '="2"'

$ws.Range("C4").Formula = '="9"'
$ws.Range("D4").Formula = '="12"'
$ws.Range("E4").Formula = '="0"'
$ws.Range("F4").Formula = '="0"'

$ws.Range("C5").Formula = '="35"'
$ws.Range("D5").Formula = '="22"'
$ws.Range("E5").Formula = '="4"'
$ws.Range("F5").Formula = '="1"'

$ws.Range("C6").Formula = '="70"'
$ws.Range("D6").Formula = '="48"'
$ws.Range("E6").Formula = '="7"'
$ws.Range("F6").Formula = '="2"'

$rng = $ws.Range("C2:F6")
$rng.Copy()
$rng.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
